$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# The "Förändrad" (Changed) date column C was bumped by one day
# (2023-09-11 -> 2023-09-12, serial 45180 -> 45181) for every data row.
$ws.Range("C2:C28").Value = 45181
